$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Cells.Item(1,26)

# Row 58
$ws.Range("A2:I2").Copy($ws.Range("A58:I58"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(58,1).PasteSpecial(-4163)
$ws.Cells.Item(58,2).Value = "08:00:00"
$ws.Cells.Item(58,3).Value = 8
$ws.Cells.Item(58,4).Value = 0
$ws.Cells.Item(58,5).Value = 1
$ws.Cells.Item(58,6).Value = 0
$ws.Cells.Item(58,7).Value = 0
$ws.Cells.Item(58,8).Value = 0
$ws.Cells.Item(58,9).Value = "Idle"

# Row 59
$ws.Range("A3:I3").Copy($ws.Range("A59:I59"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(59,1).PasteSpecial(-4163)
$ws.Cells.Item(59,2).Value = "08:00:00"
$ws.Cells.Item(59,3).Value = 8
$ws.Cells.Item(59,4).Value = 0
$ws.Cells.Item(59,5).Value = 2
$ws.Cells.Item(59,6).Value = 1
$ws.Cells.Item(59,7).Value = 0
$ws.Cells.Item(59,8).Value = 0
$ws.Cells.Item(59,9).Value = "WaitingForPassenger"

# Row 60
$ws.Range("A2:I2").Copy($ws.Range("A60:I60"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(60,1).PasteSpecial(-4163)
$ws.Cells.Item(60,2).Value = "08:00:00"
$ws.Cells.Item(60,3).Value = 8
$ws.Cells.Item(60,4).Value = 0
$ws.Cells.Item(60,5).Value = 7
$ws.Cells.Item(60,6).Value = 1
$ws.Cells.Item(60,7).Value = 0
$ws.Cells.Item(60,8).Value = 0
$ws.Cells.Item(60,9).Value = "WaitingForPassenger"

# Row 61
$ws.Range("A3:I3").Copy($ws.Range("A61:I61"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(61,1).PasteSpecial(-4163)
$ws.Cells.Item(61,2).Value = "08:00:00"
$ws.Cells.Item(61,3).Value = 8
$ws.Cells.Item(61,4).Value = 0
$ws.Cells.Item(61,5).Value = 1
$ws.Cells.Item(61,6).Value = 0
$ws.Cells.Item(61,7).Value = 3
$ws.Cells.Item(61,8).Value = 0
$ws.Cells.Item(61,9).Value = "PickedUp"

# Row 62
$ws.Range("A2:I2").Copy($ws.Range("A62:I62"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(62,1).PasteSpecial(-4163)
$ws.Cells.Item(62,2).Value = "08:00:00"
$ws.Cells.Item(62,3).Value = 8
$ws.Cells.Item(62,4).Value = 0
$ws.Cells.Item(62,5).Value = 3
$ws.Cells.Item(62,6).Value = 1
$ws.Cells.Item(62,7).Value = 0
$ws.Cells.Item(62,8).Value = 1
$ws.Cells.Item(62,9).Value = "WaitingForPassenger"

# Row 63
$ws.Range("A3:I3").Copy($ws.Range("A63:I63"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(63,1).PasteSpecial(-4163)
$ws.Cells.Item(63,2).Value = "08:00:00"
$ws.Cells.Item(63,3).Value = 8
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 2
$ws.Cells.Item(63,6).Value = 1
$ws.Cells.Item(63,7).Value = 9
$ws.Cells.Item(63,8).Value = 1
$ws.Cells.Item(63,9).Value = "PickedUp"

# Row 64
$ws.Range("A2:I2").Copy($ws.Range("A64:I64"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(64,1).PasteSpecial(-4163)
$ws.Cells.Item(64,2).Value = "08:00:00"
$ws.Cells.Item(64,3).Value = 8
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 3
$ws.Cells.Item(64,6).Value = 1
$ws.Cells.Item(64,7).Value = 12
$ws.Cells.Item(64,8).Value = 2
$ws.Cells.Item(64,9).Value = "PickedUp"

# Row 65
$ws.Range("A3:I3").Copy($ws.Range("A65:I65"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(65,1).PasteSpecial(-4163)
$ws.Cells.Item(65,2).Value = "08:00:00"
$ws.Cells.Item(65,3).Value = 8
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 7
$ws.Cells.Item(65,6).Value = 1
$ws.Cells.Item(65,7).Value = 33
$ws.Cells.Item(65,8).Value = 0
$ws.Cells.Item(65,9).Value = "PickedUp"

# Row 66
$ws.Range("A2:I2").Copy($ws.Range("A66:I66"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(66,1).PasteSpecial(-4163)
$ws.Cells.Item(66,2).Value = "08:00:00"
$ws.Cells.Item(66,3).Value = 8
$ws.Cells.Item(66,4).Value = 0
$ws.Cells.Item(66,5).Value = 14
$ws.Cells.Item(66,6).Value = 9
$ws.Cells.Item(66,7).Value = 0
$ws.Cells.Item(66,8).Value = 1
$ws.Cells.Item(66,9).Value = "MovingUp"

# Row 67
$ws.Range("A3:I3").Copy($ws.Range("A67:I67"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(67,1).PasteSpecial(-4163)
$ws.Cells.Item(67,2).Value = "08:00:00"
$ws.Cells.Item(67,3).Value = 8
$ws.Cells.Item(67,4).Value = 0
$ws.Cells.Item(67,5).Value = 14
$ws.Cells.Item(67,6).Value = 9
$ws.Cells.Item(67,7).Value = 6
$ws.Cells.Item(67,8).Value = 1
$ws.Cells.Item(67,9).Value = "PickedUp"

# Row 68
$ws.Range("A2:I2").Copy($ws.Range("A68:I68"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(68,1).PasteSpecial(-4163)
$ws.Cells.Item(68,2).Value = "08:00:00"
$ws.Cells.Item(68,3).Value = 8
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 2
$ws.Cells.Item(68,6).Value = 0
$ws.Cells.Item(68,7).Value = 0
$ws.Cells.Item(68,8).Value = 0
$ws.Cells.Item(68,9).Value = "Idle"

# Row 69
$ws.Range("A3:I3").Copy($ws.Range("A69:I69"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(69,1).PasteSpecial(-4163)
$ws.Cells.Item(69,2).Value = "08:00:00"
$ws.Cells.Item(69,3).Value = 8
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 2
$ws.Cells.Item(69,6).Value = 0
$ws.Cells.Item(69,7).Value = 4
$ws.Cells.Item(69,8).Value = 0
$ws.Cells.Item(69,9).Value = "PickedUp"

# Row 70
$ws.Range("A2:I2").Copy($ws.Range("A70:I70"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(70,1).PasteSpecial(-4163)
$ws.Cells.Item(70,2).Value = "07:58:00"
$ws.Cells.Item(70,3).Value = 7
$ws.Cells.Item(70,4).Value = 58
$ws.Cells.Item(70,5).Value = 5
$ws.Cells.Item(70,6).Value = 0
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = 0
$ws.Cells.Item(70,9).Value = "WaitingForPassenger"

# Row 71
$ws.Range("A3:I3").Copy($ws.Range("A71:I71"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(71,1).PasteSpecial(-4163)
$ws.Cells.Item(71,2).Value = "07:58:00"
$ws.Cells.Item(71,3).Value = 7
$ws.Cells.Item(71,4).Value = 58
$ws.Cells.Item(71,5).Value = 5
$ws.Cells.Item(71,6).Value = 0
$ws.Cells.Item(71,7).Value = 8
$ws.Cells.Item(71,8).Value = 0
$ws.Cells.Item(71,9).Value = "PickedUp"

# Row 72
$ws.Range("A2:I2").Copy($ws.Range("A72:I72"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(72,1).PasteSpecial(-4163)
$ws.Cells.Item(72,2).Value = "10:58:00"
$ws.Cells.Item(72,3).Value = 10
$ws.Cells.Item(72,4).Value = 58
$ws.Cells.Item(72,5).Value = 5
$ws.Cells.Item(72,6).Value = 13
$ws.Cells.Item(72,7).Value = 0
$ws.Cells.Item(72,8).Value = 0
$ws.Cells.Item(72,9).Value = "Idle"

# Row 73
$ws.Range("A3:I3").Copy($ws.Range("A73:I73"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(73,1).PasteSpecial(-4163)
$ws.Cells.Item(73,2).Value = "10:58:00"
$ws.Cells.Item(73,3).Value = 10
$ws.Cells.Item(73,4).Value = 58
$ws.Cells.Item(73,5).Value = 5
$ws.Cells.Item(73,6).Value = 13
$ws.Cells.Item(73,7).Value = 10
$ws.Cells.Item(73,8).Value = 0
$ws.Cells.Item(73,9).Value = "PickedUp"

# Row 74
$ws.Range("A2:I2").Copy($ws.Range("A74:I74"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(74,1).PasteSpecial(-4163)
$ws.Cells.Item(74,2).Value = "10:58:00"
$ws.Cells.Item(74,3).Value = 10
$ws.Cells.Item(74,4).Value = 58
$ws.Cells.Item(74,5).Value = 6
$ws.Cells.Item(74,6).Value = 5
$ws.Cells.Item(74,7).Value = 0
$ws.Cells.Item(74,8).Value = 1
$ws.Cells.Item(74,9).Value = "DoorClosing"

# Row 75
$ws.Range("A3:I3").Copy($ws.Range("A75:I75"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(75,1).PasteSpecial(-4163)
$ws.Cells.Item(75,2).Value = "10:58:00"
$ws.Cells.Item(75,3).Value = 10
$ws.Cells.Item(75,4).Value = 58
$ws.Cells.Item(75,5).Value = 6
$ws.Cells.Item(75,6).Value = 5
$ws.Cells.Item(75,7).Value = 0
$ws.Cells.Item(75,8).Value = 1
$ws.Cells.Item(75,9).Value = "DoorClosing"

# Row 76
$ws.Range("A2:I2").Copy($ws.Range("A76:I76"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(76,1).PasteSpecial(-4163)
$ws.Cells.Item(76,2).Value = "10:58:00"
$ws.Cells.Item(76,3).Value = 10
$ws.Cells.Item(76,4).Value = 58
$ws.Cells.Item(76,5).Value = 6
$ws.Cells.Item(76,6).Value = 5
$ws.Cells.Item(76,7).Value = 0
$ws.Cells.Item(76,8).Value = 1
$ws.Cells.Item(76,9).Value = "DoorClosing"

# Row 77
$ws.Range("A3:I3").Copy($ws.Range("A77:I77"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(77,1).PasteSpecial(-4163)
$ws.Cells.Item(77,2).Value = "10:58:00"
$ws.Cells.Item(77,3).Value = 10
$ws.Cells.Item(77,4).Value = 58
$ws.Cells.Item(77,5).Value = 6
$ws.Cells.Item(77,6).Value = 5
$ws.Cells.Item(77,7).Value = 0
$ws.Cells.Item(77,8).Value = 1
$ws.Cells.Item(77,9).Value = "DoorClosing"

# Row 78
$ws.Range("A2:I2").Copy($ws.Range("A78:I78"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(78,1).PasteSpecial(-4163)
$ws.Cells.Item(78,2).Value = "05:52:00"
$ws.Cells.Item(78,3).Value = 5
$ws.Cells.Item(78,4).Value = 52
$ws.Cells.Item(78,5).Value = 3
$ws.Cells.Item(78,6).Value = 0
$ws.Cells.Item(78,7).Value = 0
$ws.Cells.Item(78,8).Value = 0
$ws.Cells.Item(78,9).Value = "Idle"

# Row 79
$ws.Range("A3:I3").Copy($ws.Range("A79:I79"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(79,1).PasteSpecial(-4163)
$ws.Cells.Item(79,2).Value = "05:52:00"
$ws.Cells.Item(79,3).Value = 5
$ws.Cells.Item(79,4).Value = 52
$ws.Cells.Item(79,5).Value = 5
$ws.Cells.Item(79,6).Value = 1
$ws.Cells.Item(79,7).Value = 0
$ws.Cells.Item(79,8).Value = 0
$ws.Cells.Item(79,9).Value = "MovingUp"

# Row 80
$ws.Range("A2:I2").Copy($ws.Range("A80:I80"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(80,1).PasteSpecial(-4163)
$ws.Cells.Item(80,2).Value = "05:52:00"
$ws.Cells.Item(80,3).Value = 5
$ws.Cells.Item(80,4).Value = 52
$ws.Cells.Item(80,5).Value = 4
$ws.Cells.Item(80,6).Value = 2
$ws.Cells.Item(80,7).Value = 0
$ws.Cells.Item(80,8).Value = 0
$ws.Cells.Item(80,9).Value = "MovingUp"

# Row 81
$ws.Range("A3:I3").Copy($ws.Range("A81:I81"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(81,1).PasteSpecial(-4163)
$ws.Cells.Item(81,2).Value = "05:52:00"
$ws.Cells.Item(81,3).Value = 5
$ws.Cells.Item(81,4).Value = 52
$ws.Cells.Item(81,5).Value = 3
$ws.Cells.Item(81,6).Value = 0
$ws.Cells.Item(81,7).Value = 5
$ws.Cells.Item(81,8).Value = 0
$ws.Cells.Item(81,9).Value = "PickedUp"

# Row 82
$ws.Range("A2:I2").Copy($ws.Range("A82:I82"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(82,1).PasteSpecial(-4163)
$ws.Cells.Item(82,2).Value = "05:52:00"
$ws.Cells.Item(82,3).Value = 5
$ws.Cells.Item(82,4).Value = 52
$ws.Cells.Item(82,5).Value = 4
$ws.Cells.Item(82,6).Value = 2
$ws.Cells.Item(82,7).Value = 13
$ws.Cells.Item(82,8).Value = 1
$ws.Cells.Item(82,9).Value = "PickedUp"

# Row 83
$ws.Range("A3:I3").Copy($ws.Range("A83:I83"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(83,1).PasteSpecial(-4163)
$ws.Cells.Item(83,2).Value = "05:52:00"
$ws.Cells.Item(83,3).Value = 5
$ws.Cells.Item(83,4).Value = 52
$ws.Cells.Item(83,5).Value = 5
$ws.Cells.Item(83,6).Value = 1
$ws.Cells.Item(83,7).Value = 21
$ws.Cells.Item(83,8).Value = 2
$ws.Cells.Item(83,9).Value = "PickedUp"

# Row 84
$ws.Range("A2:I2").Copy($ws.Range("A84:I84"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(84,1).PasteSpecial(-4163)
$ws.Cells.Item(84,2).Value = "08:00:00"
$ws.Cells.Item(84,3).Value = 8
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 2
$ws.Cells.Item(84,6).Value = 0
$ws.Cells.Item(84,7).Value = 0
$ws.Cells.Item(84,8).Value = 0
$ws.Cells.Item(84,9).Value = "Idle"

# Row 85
$ws.Range("A3:I3").Copy($ws.Range("A85:I85"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(85,1).PasteSpecial(-4163)
$ws.Cells.Item(85,2).Value = "08:00:00"
$ws.Cells.Item(85,3).Value = 8
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 2
$ws.Cells.Item(85,6).Value = 0
$ws.Cells.Item(85,7).Value = 4
$ws.Cells.Item(85,8).Value = 0
$ws.Cells.Item(85,9).Value = "PickedUp"

# Row 86
$ws.Range("A2:I2").Copy($ws.Range("A86:I86"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(86,1).PasteSpecial(-4163)
$ws.Cells.Item(86,2).Value = "08:00:00"
$ws.Cells.Item(86,3).Value = 8
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 3
$ws.Cells.Item(86,6).Value = 0
$ws.Cells.Item(86,7).Value = 0
$ws.Cells.Item(86,8).Value = 0
$ws.Cells.Item(86,9).Value = "Idle"

# Row 87
$ws.Range("A3:I3").Copy($ws.Range("A87:I87"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(87,1).PasteSpecial(-4163)
$ws.Cells.Item(87,2).Value = "08:00:00"
$ws.Cells.Item(87,3).Value = 8
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 3
$ws.Cells.Item(87,6).Value = 0
$ws.Cells.Item(87,7).Value = 6
$ws.Cells.Item(87,8).Value = 0
$ws.Cells.Item(87,9).Value = "PickedUp"

# Row 88
$ws.Range("A2:I2").Copy($ws.Range("A88:I88"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(88,1).PasteSpecial(-4163)
$ws.Cells.Item(88,2).Value = "08:00:00"
$ws.Cells.Item(88,3).Value = 8
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 3
$ws.Cells.Item(88,6).Value = 0
$ws.Cells.Item(88,7).Value = 0
$ws.Cells.Item(88,8).Value = 0
$ws.Cells.Item(88,9).Value = "Idle"

# Row 89
$ws.Range("A3:I3").Copy($ws.Range("A89:I89"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(89,1).PasteSpecial(-4163)
$ws.Cells.Item(89,2).Value = "08:00:00"
$ws.Cells.Item(89,3).Value = 8
$ws.Cells.Item(89,4).Value = 0
$ws.Cells.Item(89,5).Value = 3
$ws.Cells.Item(89,6).Value = 0
$ws.Cells.Item(89,7).Value = 5
$ws.Cells.Item(89,8).Value = 0
$ws.Cells.Item(89,9).Value = "PickedUp"

# Row 90
$ws.Range("A2:I2").Copy($ws.Range("A90:I90"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(90,1).PasteSpecial(-4163)
$ws.Cells.Item(90,2).Value = "08:00:00"
$ws.Cells.Item(90,3).Value = 8
$ws.Cells.Item(90,4).Value = 0
$ws.Cells.Item(90,5).Value = 1
$ws.Cells.Item(90,6).Value = 0
$ws.Cells.Item(90,7).Value = 0
$ws.Cells.Item(90,8).Value = 0
$ws.Cells.Item(90,9).Value = "Idle"

# Row 91
$ws.Range("A3:I3").Copy($ws.Range("A91:I91"))
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(91,1).PasteSpecial(-4163)
$ws.Cells.Item(91,2).Value = "08:00:00"
$ws.Cells.Item(91,3).Value = 8
$ws.Cells.Item(91,4).Value = 0
$ws.Cells.Item(91,5).Value = 1
$ws.Cells.Item(91,6).Value = 0
$ws.Cells.Item(91,7).Value = 3
$ws.Cells.Item(91,8).Value = 0
$ws.Cells.Item(91,9).Value = "PickedUp"

# Row 92
$ws.Range("A2:I2").Copy($ws.Range("A92:I92"))
$scratch.NumberFormat = "@"
$scratch.Value = "2026-01-18"
$scratch.Copy()
$ws.Cells.Item(92,1).PasteSpecial(-4163)
$ws.Cells.Item(92,2).Value = "08:00:00"
$ws.Cells.Item(92,3).Value = 8
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 3
$ws.Cells.Item(92,6).Value = 0
$ws.Cells.Item(92,7).Value = 0
$ws.Cells.Item(92,8).Value = 0
$ws.Cells.Item(92,9).Value = "Idle"

# Row 93
$ws.Range("A3:I3").Copy($ws.Range("A93:I93"))
$scratch.Copy()
$ws.Cells.Item(93,1).PasteSpecial(-4163)
$ws.Cells.Item(93,2).Value = "08:00:00"
$ws.Cells.Item(93,3).Value = 8
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 3
$ws.Cells.Item(93,6).Value = 0
$ws.Cells.Item(93,7).Value = 5
$ws.Cells.Item(93,8).Value = 0
$ws.Cells.Item(93,9).Value = "PickedUp"

# Row 94
$ws.Range("A2:I2").Copy($ws.Range("A94:I94"))
$scratch.Copy()
$ws.Cells.Item(94,1).PasteSpecial(-4163)
$ws.Cells.Item(94,2).Value = "08:00:00"
$ws.Cells.Item(94,3).Value = 8
$ws.Cells.Item(94,4).Value = 0
$ws.Cells.Item(94,5).Value = 2
$ws.Cells.Item(94,6).Value = 0
$ws.Cells.Item(94,7).Value = 0
$ws.Cells.Item(94,8).Value = 0
$ws.Cells.Item(94,9).Value = "Idle"

# Row 95
$ws.Range("A3:I3").Copy($ws.Range("A95:I95"))
$scratch.Copy()
$ws.Cells.Item(95,1).PasteSpecial(-4163)
$ws.Cells.Item(95,2).Value = "08:00:00"
$ws.Cells.Item(95,3).Value = 8
$ws.Cells.Item(95,4).Value = 0
$ws.Cells.Item(95,5).Value = 3
$ws.Cells.Item(95,6).Value = 2
$ws.Cells.Item(95,7).Value = 0
$ws.Cells.Item(95,8).Value = 0
$ws.Cells.Item(95,9).Value = "WaitingForPassenger"

# Row 96
$ws.Range("A2:I2").Copy($ws.Range("A96:I96"))
$scratch.Copy()
$ws.Cells.Item(96,1).PasteSpecial(-4163)
$ws.Cells.Item(96,2).Value = "08:00:00"
$ws.Cells.Item(96,3).Value = 8
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 4
$ws.Cells.Item(96,6).Value = 2
$ws.Cells.Item(96,7).Value = 0
$ws.Cells.Item(96,8).Value = 0
$ws.Cells.Item(96,9).Value = "WaitingForPassenger"

# Row 97
$ws.Range("A3:I3").Copy($ws.Range("A97:I97"))
$scratch.Copy()
$ws.Cells.Item(97,1).PasteSpecial(-4163)
$ws.Cells.Item(97,2).Value = "08:00:00"
$ws.Cells.Item(97,3).Value = 8
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 5
$ws.Cells.Item(97,6).Value = 2
$ws.Cells.Item(97,7).Value = 0
$ws.Cells.Item(97,8).Value = 0
$ws.Cells.Item(97,9).Value = "WaitingForPassenger"

# Row 98
$ws.Range("A2:I2").Copy($ws.Range("A98:I98"))
$scratch.Copy()
$ws.Cells.Item(98,1).PasteSpecial(-4163)
$ws.Cells.Item(98,2).Value = "08:00:00"
$ws.Cells.Item(98,3).Value = 8
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 5
$ws.Cells.Item(98,6).Value = 2
$ws.Cells.Item(98,7).Value = 0
$ws.Cells.Item(98,8).Value = 0
$ws.Cells.Item(98,9).Value = "WaitingForPassenger"

# Row 99
$ws.Range("A3:I3").Copy($ws.Range("A99:I99"))
$scratch.Copy()
$ws.Cells.Item(99,1).PasteSpecial(-4163)
$ws.Cells.Item(99,2).Value = "08:00:00"
$ws.Cells.Item(99,3).Value = 8
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 2
$ws.Cells.Item(99,6).Value = 0
$ws.Cells.Item(99,7).Value = 4
$ws.Cells.Item(99,8).Value = 0
$ws.Cells.Item(99,9).Value = "PickedUp"

# Row 100
$ws.Range("A2:I2").Copy($ws.Range("A100:I100"))
$scratch.Copy()
$ws.Cells.Item(100,1).PasteSpecial(-4163)
$ws.Cells.Item(100,2).Value = "08:00:00"
$ws.Cells.Item(100,3).Value = 8
$ws.Cells.Item(100,4).Value = 0
$ws.Cells.Item(100,5).Value = 3
$ws.Cells.Item(100,6).Value = 2
$ws.Cells.Item(100,7).Value = 9
$ws.Cells.Item(100,8).Value = 1
$ws.Cells.Item(100,9).Value = "PickedUp"

# Row 101
$ws.Range("A3:I3").Copy($ws.Range("A101:I101"))
$scratch.Copy()
$ws.Cells.Item(101,1).PasteSpecial(-4163)
$ws.Cells.Item(101,2).Value = "08:00:00"
$ws.Cells.Item(101,3).Value = 8
$ws.Cells.Item(101,4).Value = 0
$ws.Cells.Item(101,5).Value = 4
$ws.Cells.Item(101,6).Value = 2
$ws.Cells.Item(101,7).Value = 16
$ws.Cells.Item(101,8).Value = 2
$ws.Cells.Item(101,9).Value = "PickedUp"

# Row 102
$ws.Range("A2:I2").Copy($ws.Range("A102:I102"))
$scratch.Copy()
$ws.Cells.Item(102,1).PasteSpecial(-4163)
$ws.Cells.Item(102,2).Value = "08:00:00"
$ws.Cells.Item(102,3).Value = 8
$ws.Cells.Item(102,4).Value = 0
$ws.Cells.Item(102,5).Value = 5
$ws.Cells.Item(102,6).Value = 2
$ws.Cells.Item(102,7).Value = 48
$ws.Cells.Item(102,8).Value = 2
$ws.Cells.Item(102,9).Value = "PickedUp"

# Row 103
$ws.Range("A3:I3").Copy($ws.Range("A103:I103"))
$scratch.Copy()
$ws.Cells.Item(103,1).PasteSpecial(-4163)
$ws.Cells.Item(103,2).Value = "08:00:00"
$ws.Cells.Item(103,3).Value = 8
$ws.Cells.Item(103,4).Value = 0
$ws.Cells.Item(103,5).Value = 5
$ws.Cells.Item(103,6).Value = 2
$ws.Cells.Item(103,7).Value = 55
$ws.Cells.Item(103,8).Value = 3
$ws.Cells.Item(103,9).Value = "PickedUp"

# Row 104
$ws.Range("A2:I2").Copy($ws.Range("A104:I104"))
$scratch.Copy()
$ws.Cells.Item(104,1).PasteSpecial(-4163)
$ws.Cells.Item(104,2).Value = "08:00:00"
$ws.Cells.Item(104,3).Value = 8
$ws.Cells.Item(104,4).Value = 0
$ws.Cells.Item(104,5).Value = 2
$ws.Cells.Item(104,6).Value = 0
$ws.Cells.Item(104,7).Value = 0
$ws.Cells.Item(104,8).Value = 0
$ws.Cells.Item(104,9).Value = "Idle"

# Row 105
$ws.Range("A3:I3").Copy($ws.Range("A105:I105"))
$scratch.Copy()
$ws.Cells.Item(105,1).PasteSpecial(-4163)
$ws.Cells.Item(105,2).Value = "08:00:00"
$ws.Cells.Item(105,3).Value = 8
$ws.Cells.Item(105,4).Value = 0
$ws.Cells.Item(105,5).Value = 4
$ws.Cells.Item(105,6).Value = 2
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = 0
$ws.Cells.Item(105,9).Value = "DoorOpening"

# Row 106
$ws.Range("A2:I2").Copy($ws.Range("A106:I106"))
$scratch.Copy()
$ws.Cells.Item(106,1).PasteSpecial(-4163)
$ws.Cells.Item(106,2).Value = "08:00:00"
$ws.Cells.Item(106,3).Value = 8
$ws.Cells.Item(106,4).Value = 0
$ws.Cells.Item(106,5).Value = 2
$ws.Cells.Item(106,6).Value = 0
$ws.Cells.Item(106,7).Value = 4
$ws.Cells.Item(106,8).Value = 0
$ws.Cells.Item(106,9).Value = "PickedUp"

# Row 107
$ws.Range("A3:I3").Copy($ws.Range("A107:I107"))
$scratch.Copy()
$ws.Cells.Item(107,1).PasteSpecial(-4163)
$ws.Cells.Item(107,2).Value = "08:00:00"
$ws.Cells.Item(107,3).Value = 8
$ws.Cells.Item(107,4).Value = 0
$ws.Cells.Item(107,5).Value = 4
$ws.Cells.Item(107,6).Value = 0
$ws.Cells.Item(107,7).Value = 0
$ws.Cells.Item(107,8).Value = 0
$ws.Cells.Item(107,9).Value = "Idle"

# Row 108
$ws.Range("A2:I2").Copy($ws.Range("A108:I108"))
$scratch.Copy()
$ws.Cells.Item(108,1).PasteSpecial(-4163)
$ws.Cells.Item(108,2).Value = "08:00:00"
$ws.Cells.Item(108,3).Value = 8
$ws.Cells.Item(108,4).Value = 0
$ws.Cells.Item(108,5).Value = 4
$ws.Cells.Item(108,6).Value = 0
$ws.Cells.Item(108,7).Value = 6
$ws.Cells.Item(108,8).Value = 0
$ws.Cells.Item(108,9).Value = "PickedUp"

# Row 109
$ws.Range("A3:I3").Copy($ws.Range("A109:I109"))
$scratch.Copy()
$ws.Cells.Item(109,1).PasteSpecial(-4163)
$ws.Cells.Item(109,2).Value = "08:00:00"
$ws.Cells.Item(109,3).Value = 8
$ws.Cells.Item(109,4).Value = 0
$ws.Cells.Item(109,5).Value = 2
$ws.Cells.Item(109,6).Value = 0
$ws.Cells.Item(109,7).Value = 0
$ws.Cells.Item(109,8).Value = 0
$ws.Cells.Item(109,9).Value = "Idle"

# Row 110
$ws.Range("A2:I2").Copy($ws.Range("A110:I110"))
$scratch.Copy()
$ws.Cells.Item(110,1).PasteSpecial(-4163)
$ws.Cells.Item(110,2).Value = "08:00:00"
$ws.Cells.Item(110,3).Value = 8
$ws.Cells.Item(110,4).Value = 0
$ws.Cells.Item(110,5).Value = 2
$ws.Cells.Item(110,6).Value = 0
$ws.Cells.Item(110,7).Value = 4
$ws.Cells.Item(110,8).Value = 0
$ws.Cells.Item(110,9).Value = "PickedUp"

# Row 111
$ws.Range("A3:I3").Copy($ws.Range("A111:I111"))
$scratch.Copy()
$ws.Cells.Item(111,1).PasteSpecial(-4163)
$ws.Cells.Item(111,2).Value = "23:45:00"
$ws.Cells.Item(111,3).Value = 23
$ws.Cells.Item(111,4).Value = 45
$ws.Cells.Item(111,5).Value = 0
$ws.Cells.Item(111,6).Value = 0
$ws.Cells.Item(111,7).Value = 0
$ws.Cells.Item(111,8).Value = 0
$ws.Cells.Item(111,9).Value = "Idle"

# Row 112
$ws.Range("A2:I2").Copy($ws.Range("A112:I112"))
$scratch.Copy()
$ws.Cells.Item(112,1).PasteSpecial(-4163)
$ws.Cells.Item(112,2).Value = "23:45:00"
$ws.Cells.Item(112,3).Value = 23
$ws.Cells.Item(112,4).Value = 45
$ws.Cells.Item(112,5).Value = 0
$ws.Cells.Item(112,6).Value = 0
$ws.Cells.Item(112,7).Value = 2
$ws.Cells.Item(112,8).Value = 0
$ws.Cells.Item(112,9).Value = "PickedUp"

# Row 113
$ws.Range("A3:I3").Copy($ws.Range("A113:I113"))
$scratch.Copy()
$ws.Cells.Item(113,1).PasteSpecial(-4163)
$ws.Cells.Item(113,2).Value = "23:45:00"
$ws.Cells.Item(113,3).Value = 23
$ws.Cells.Item(113,4).Value = 45
$ws.Cells.Item(113,5).Value = 2
$ws.Cells.Item(113,6).Value = 1
$ws.Cells.Item(113,7).Value = 0
$ws.Cells.Item(113,8).Value = 0
$ws.Cells.Item(113,9).Value = "Idle"

# Row 114
$ws.Range("A2:I2").Copy($ws.Range("A114:I114"))
$scratch.Copy()
$ws.Cells.Item(114,1).PasteSpecial(-4163)
$ws.Cells.Item(114,2).Value = "23:45:00"
$ws.Cells.Item(114,3).Value = 23
$ws.Cells.Item(114,4).Value = 45
$ws.Cells.Item(114,5).Value = 2
$ws.Cells.Item(114,6).Value = 1
$ws.Cells.Item(114,7).Value = 3
$ws.Cells.Item(114,8).Value = 0
$ws.Cells.Item(114,9).Value = "PickedUp"

# Row 115
$ws.Range("A3:I3").Copy($ws.Range("A115:I115"))
$scratch.Copy()
$ws.Cells.Item(115,1).PasteSpecial(-4163)
$ws.Cells.Item(115,2).Value = "23:45:00"
$ws.Cells.Item(115,3).Value = 23
$ws.Cells.Item(115,4).Value = 45
$ws.Cells.Item(115,5).Value = 19
$ws.Cells.Item(115,6).Value = 5
$ws.Cells.Item(115,7).Value = 0
$ws.Cells.Item(115,8).Value = 0
$ws.Cells.Item(115,9).Value = "Idle"

# Row 116
$ws.Range("A2:I2").Copy($ws.Range("A116:I116"))
$scratch.Copy()
$ws.Cells.Item(116,1).PasteSpecial(-4163)
$ws.Cells.Item(116,2).Value = "23:45:00"
$ws.Cells.Item(116,3).Value = 23
$ws.Cells.Item(116,4).Value = 45
$ws.Cells.Item(116,5).Value = 19
$ws.Cells.Item(116,6).Value = 5
$ws.Cells.Item(116,7).Value = 16
$ws.Cells.Item(116,8).Value = 0
$ws.Cells.Item(116,9).Value = "PickedUp"

# Clean up scratch cell so it leaves no residue in the used range
$scratch.Clear()
